{"js": "// Replace the 25 two-digit-divided-by-one-digit division problems in the\n// practice table with a freshly generated set of problems, cell by cell,\n// preserving each cell's existing paragraph/run formatting.\n\n// (row, col, newText) - table has 20 rows x 5 cols; only rows 0,4,8,12,16\n// contain the practice problems (one \"block\" of 5 problems per used row).\nconst replacements = [\n  [0, 0, \"96\u00f78=12, 0\"],\n  [0, 1, \"90\u00f77=12, 6\"],\n  [0, 2, \"81\u00f78=10, 1\"],\n  [0, 3, \"54\u00f74=13, 2\"],\n  [0, 4, \"54\u00f77=7, 5\"],\n\n  [4, 0, \"32\u00f79=3, 5\"],\n  [4, 1, \"38\u00f78=4, 6\"],\n  [4, 2, \"88\u00f76=14, 4\"],\n  [4, 3, \"26\u00f76=4, 2\"],\n  [4, 4, \"40\u00f75=8, 0\"],\n\n  [8, 0, \"18\u00f79=2, 0\"],\n  [8, 1, \"26\u00f72=13, 0\"],\n  [8, 2, \"17\u00f78=2, 1\"],\n  [8, 3, \"74\u00f75=14, 4\"],\n  [8, 4, \"60\u00f74=15, 0\"],\n\n  [12, 0, \"57\u00f77=8, 1\"],\n  [12, 1, \"91\u00f79=10, 1\"],\n  [12, 2, \"26\u00f79=2, 8\"],\n  [12, 3, \"52\u00f79=5, 7\"],\n  [12, 4, \"41\u00f75=8, 1\"],\n\n  [16, 0, \"72\u00f75=14, 2\"],\n  [16, 1, \"60\u00f73=20, 0\"],\n  [16, 2, \"91\u00f78=11, 3\"],\n  [16, 3, \"26\u00f78=3, 2\"],\n  [16, 4, \"17\u00f77=2, 3\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Grab every target cell's first paragraph range up-front, then sync once\n// before mutating so we don't repeatedly round-trip.\nconst ranges = replacements.map(([row, col]) => {\n  const cell = table.getCell(row, col);\n  const para = cell.body.paragraphs.getFirst();\n  return para.getRange();\n});\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [, , newText] = replacements[i];\n  ranges[i].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 two-digit-divided-by-one-digit division problems in the\n# practice table with a freshly generated set of problems, cell by cell,\n# preserving each cell's existing paragraph/run formatting (setting\n# Range.Text in place keeps the surrounding w:rPr/w:pPr untouched).\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# (row, col, newText) - table has 20 rows x 5 cols (1-based here); only\n# rows 1,5,9,13,17 contain the practice problems (one \"block\" of 5 per row).\n$replacements = @(\n    @(1, 1, \"96\u00f78=12, 0\"),\n    @(1, 2, \"90\u00f77=12, 6\"),\n    @(1, 3, \"81\u00f78=10, 1\"),\n    @(1, 4, \"54\u00f74=13, 2\"),\n    @(1, 5, \"54\u00f77=7, 5\"),\n\n    @(5, 1, \"32\u00f79=3, 5\"),\n    @(5, 2, \"38\u00f78=4, 6\"),\n    @(5, 3, \"88\u00f76=14, 4\"),\n    @(5, 4, \"26\u00f76=4, 2\"),\n    @(5, 5, \"40\u00f75=8, 0\"),\n\n    @(9, 1, \"18\u00f79=2, 0\"),\n    @(9, 2, \"26\u00f72=13, 0\"),\n    @(9, 3, \"17\u00f78=2, 1\"),\n    @(9, 4, \"74\u00f75=14, 4\"),\n    @(9, 5, \"60\u00f74=15, 0\"),\n\n    @(13, 1, \"57\u00f77=8, 1\"),\n    @(13, 2, \"91\u00f79=10, 1\"),\n    @(13, 3, \"26\u00f79=2, 8\"),\n    @(13, 4, \"52\u00f79=5, 7\"),\n    @(13, 5, \"41\u00f75=8, 1\"),\n\n    @(17, 1, \"72\u00f75=14, 2\"),\n    @(17, 2, \"60\u00f73=20, 0\"),\n    @(17, 3, \"91\u00f78=11, 3\"),\n    @(17, 4, \"26\u00f78=3, 2\"),\n    @(17, 5, \"17\u00f77=2, 3\")\n)\n\nforeach ($item in $replacements) {\n    $row = $item[0]\n    $col = $item[1]\n    $newText = $item[2]\n    $cell = $t.Cell($row, $col)\n    $cell.Range.Text = $newText\n}\n"}
